$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix typo in the "Fotos de la implementación!" label (row 6)
$ws.Range("B6").Value = "Fotos de la implementación!"

# Mark the four measurement rows as "Completado"
$ws.Range("E3").Value = "Completado"
$ws.Range("E4").Value = "Completado"
$ws.Range("E5").Value = "Completado"
$ws.Range("E6").Value = "Completado"

# Add the new measurement notes/comments
$ws.Range("F5").Value = "Para el caso sincrónico. La amarilla es la entrada del clock, la verde es Q2, la azul es Q1 y la rosa es Q0. Se tuvo que medir usando Single. Se detecto que las salidas de los flip flops eran de 4V."
$ws.Range("F4").Value = "Mismo caso que los otros del sincrónico, las salidas no llegan a 5 V. Amarilla es clock, verde es Q2, azul es Q1 y luego rosa es Q0."

# Resize rows to fit the newly added wrapped text
$ws.Rows.Item(4).RowHeight = 45
$ws.Rows.Item(5).RowHeight = 75

# Move the active selection to E6
$ws.Range("E6").Select() | Out-Null
